# daily auto push: 2026-01-20 13:53 UTC
#
# Two new 2026/01/20 readings were appended to that day's block (which used
# to end at row 689). Everything from the old row 690 onward shifts down by
# two rows to make room, and the sheet grows from A1:D731 to A1:D733.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last existing 2026/01/20 row (689) twice, inserting the
# copies right below it. Copy+Insert (rather than a plain blank-row insert
# followed by setting .Value) keeps the date/weekday cells as plain text,
# matching the rest of the column instead of Excel reinterpreting the
# "2026/01/20" string as a date serial number.
$ws.Rows.Item(689).Copy()
$ws.Rows.Item(690).Insert()

$ws.Rows.Item(689).Copy()
$ws.Rows.Item(691).Insert()

# Fix up the two new rows' time / ranking figures.
$ws.Cells.Item(690, 3).Value = 18
$ws.Cells.Item(690, 4).Value = 174

$ws.Cells.Item(691, 3).Value = 19
$ws.Cells.Item(691, 4).Value = 172
